# RPA datasets push 2024-06-18
# Insert the newest IPO record (라메디텍) as the new second row (row 2),
# pushing the existing data rows down by one, and drop the oldest record
# that falls off the bottom of the tracked window (old row 14: 하나33호스팩).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (2-14) down by inserting a new blank row at row 2.
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherits the header row's bold/border/centered
# style; strip that back to the plain (unstyled) look used by every other
# data row before writing values into it.
$ws.Range("A2:T2").ClearFormats()

# Populate the new row 2 with the latest IPO entry. The date-looking
# columns (A, D, E) are stored as plain text in this sheet, so force a
# text number format before assigning to avoid Excel auto-converting them
# into date serial numbers.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "2024-06-05"
$ws.Range("B2").Value = "라메디텍"
$ws.Range("C2").Value = "대신"
$ws.Range("D2").Value = "2024-06-11"
$ws.Range("E2").Value = "2024-06-17"
$ws.Range("F2").Value = 20768000
$ws.Range("G2").Value = 1298000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 10400
$ws.Range("J2").Value = 12700
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "2140 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"

# The table keeps a fixed window of rows, so the oldest record (now pushed
# down to row 15, originally row 14: 하나33호스팩) is removed entirely.
$ws.Rows.Item(15).Delete()
